# Controle de Ponto - append new time-clock punch rows (160-171)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 160 only has columns A and B filled in (a stray / malformed punch entry,
# note the day value also carries the time and trailing space).
$ws.Cells.Item(160, 1).Value = "15 /05/ 2023 17:02:55 "
$ws.Cells.Item(160, 2).Value = "15 /05/ 2023 17:02:54 "

# Remaining rows (161-171) are full entry/exit pairs with a computed duration
# (in fraction-of-a-day units, formatted as [hh]:mm:ss) in column E.
# NB: scientific notation literals (1.23e-05) are written out in plain
# decimal form below because the PS parser here does not accept exponents.
$data = @(
    @{ Row = 161; A = "15 /05/ 2023"; B = "17:09:25"; C = "15 /05/ 2023"; D = "17:09:26"; E = 0.00001157407407407407 },
    @{ Row = 162; A = "15 /05/ 2023"; B = "17:09:27"; C = "15 /05/ 2023"; D = "17:09:27"; E = 0.0 },
    @{ Row = 163; A = "15 /05/ 2023"; B = "17:09:29"; C = "15 /05/ 2023"; D = "17:09:29"; E = 0.0 },
    @{ Row = 164; A = "15 /05/ 2023"; B = "17:09:31"; C = "15 /05/ 2023"; D = "17:09:29"; E = -0.00002314814814814815 },
    @{ Row = 165; A = "15 /05/ 2023"; B = "17:09:32"; C = "15 /05/ 2023"; D = "17:09:30"; E = -0.00002314814814814815 },
    @{ Row = 166; A = "15 /05/ 2023"; B = "17:09:25"; C = "15 /05/ 2023"; D = "17:09:26"; E = 0.00001157407407407407 },
    @{ Row = 167; A = "15 /05/ 2023"; B = "17:09:27"; C = "15 /05/ 2023"; D = "17:09:27"; E = 0.0 },
    @{ Row = 168; A = "15 /05/ 2023"; B = "17:09:29"; C = "15 /05/ 2023"; D = "17:09:29"; E = 0.0 },
    @{ Row = 169; A = "15 /05/ 2023"; B = "17:09:31"; C = "15 /05/ 2023"; D = "17:09:29"; E = -0.00002314814814814815 },
    @{ Row = 170; A = "15 /05/ 2023"; B = "17:09:32"; C = "15 /05/ 2023"; D = "17:09:30"; E = -0.00002314814814814815 },
    @{ Row = 171; A = "15 /05/ 2023"; B = "17:09:46"; C = "15 /05/ 2023"; D = "17:09:30"; E = 0.00001157407407407407 }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 1).Value = $entry.A
    $ws.Cells.Item($r, 2).Value = $entry.B
    $ws.Cells.Item($r, 3).Value = $entry.C
    $ws.Cells.Item($r, 4).Value = $entry.D
    $eCell = $ws.Cells.Item($r, 5)
    $eCell.Value = $entry.E
    $eCell.NumberFormat = "[hh]:mm:ss"
}
